# Weekly update: two new price records were added to the
# "Hortaliza, Vega Central Mapocho de Santiago - Zapallo italiano" sheet.
# Inserting two rows at 623/624 pushes the existing rows 623-638 down to
# 625-640 (their content is unchanged), and the two freshly inserted rows
# are populated with the new weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 623 (this shifts every
# row from 623 downward by two, matching rows 625-640 in the target file).
$ws.Range("A623:A624").EntireRow.Insert()

# --- New row 623 ---------------------------------------------------------
$ws.Cells.Item(623, 1).Value = 9
$ws.Cells.Item(623, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(623, 3).Value = "Metropolitana"
$ws.Cells.Item(623, 4).Value2 = 45239
$ws.Cells.Item(623, 5).Value = 13
$ws.Cells.Item(623, 6).Value = 100112032
$ws.Cells.Item(623, 7).Value = "Zapallo italiano"
$ws.Cells.Item(623, 8).Value = "Bola 8"
$ws.Cells.Item(623, 9).Value = "Primera"
$ws.Cells.Item(623, 10).Value = 70
$ws.Cells.Item(623, 11).Value = 13000
$ws.Cells.Item(623, 12).Value = 14000
$ws.Cells.Item(623, 13).Value = 13500
$ws.Cells.Item(623, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(623, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(623, 16).Value = 270
$ws.Cells.Item(623, 17).Value = 50
$ws.Cells.Item(623, 18).Value = "Hortaliza"

# --- New row 624 ---------------------------------------------------------
$ws.Cells.Item(624, 1).Value = 9
$ws.Cells.Item(624, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(624, 3).Value = "Metropolitana"
$ws.Cells.Item(624, 4).Value2 = 45239
$ws.Cells.Item(624, 5).Value = 13
$ws.Cells.Item(624, 6).Value = 100112032
$ws.Cells.Item(624, 7).Value = "Zapallo italiano"
$ws.Cells.Item(624, 8).Value = "Sin especificar"
$ws.Cells.Item(624, 9).Value = "Primera"
$ws.Cells.Item(624, 10).Value = 160
$ws.Cells.Item(624, 11).Value = 17000
$ws.Cells.Item(624, 12).Value = 18000
$ws.Cells.Item(624, 13).Value = 17500
$ws.Cells.Item(624, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(624, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(624, 16).Value = 350
$ws.Cells.Item(624, 17).Value = 50
$ws.Cells.Item(624, 18).Value = "Hortaliza"
